$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$modelText = "MultiOutputRegressor(estimator=GridSearchCV(cv=5,
                                            estimator=Pipeline(steps=[('model',
                                                                       RandomForestRegressor())]),
                                            param_grid={'model__max_depth': [3,
                                                                             5,
                                                                             7],
                                                        'model__n_estimators': [50,
                                                                                100,
                                                                                150]},
                                            scoring='neg_mean_squared_error'))"

# Header cell F1
$ws.Range("F1").Value = "Modelo"
$ws.Range("A1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Update numeric values for B2:D5
$ws.Range("B2").Value = 0.5007975339915302
$ws.Range("C2").Value = 0.9900280575939385
$ws.Range("D2").Value = 0.577282691421273

$ws.Range("B3").Value = 0.2471391001939683
$ws.Range("C3").Value = 0.9951696600716268
$ws.Range("D3").Value = 0.3900443243076906

$ws.Range("B4").Value = 0.2842819667707182
$ws.Range("C4").Value = 0.9945310471005109
$ws.Range("D4").Value = 0.4325455182881846

$ws.Range("B5").Value = 0.4115775457434243
$ws.Range("C5").Value = 0.9918842542978411
$ws.Range("D5").Value = 0.4928213322493636

# Model descriptions in column F for rows 2-5
$ws.Range("F2").Value = $modelText
$ws.Range("F3").Value = $modelText
$ws.Range("F4").Value = $modelText
$ws.Range("F5").Value = $modelText

# Restore natural row heights (undo auto-height expansion caused by
# the embedded newlines in the long model description text)
$ws.Rows.Item(2).AutoFit()
$ws.Rows.Item(3).AutoFit()
$ws.Rows.Item(4).AutoFit()
$ws.Rows.Item(5).AutoFit()
